$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("C2").Value = 28
$ws.Range("K2").Value = 69
$ws.Range("L2").Value = 67
$ws.Range("J3").Value = 84
$ws.Range("F4").Value = 3
$ws.Range("C6").Value = 223
$ws.Range("D6").Value = 199
$ws.Range("E6").Value = 208
$ws.Range("G6").Value = 223
$ws.Range("H6").Value = 198
$ws.Range("I6").Value = 263
$ws.Range("J6").Value = 192
$ws.Range("K6").Value = 228
$ws.Range("L6").Value = 238
$ws.Range("C7").Value = 302
$ws.Range("D7").Value = 313
$ws.Range("E7").Value = 312
$ws.Range("F7").Value = 342
$ws.Range("G7").Value = 329
$ws.Range("H7").Value = 304
$ws.Range("I7").Value = 407
$ws.Range("J7").Value = 349
$ws.Range("K7").Value = 414
$ws.Range("L7").Value = 424

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("C7").Value = 24
$ws.Range("J7").Value = 25
$ws.Range("L7").Value = 14
$ws.Range("D25").Value = 4
$ws.Range("D28").Value = 4
$ws.Range("K28").Value = 8
$ws.Range("F30").Value = 29
$ws.Range("J30").Value = 15
$ws.Range("C34").Value = 15
$ws.Range("E34").Value = 13
$ws.Range("I34").Value = 22
$ws.Range("K34").Value = 29
$ws.Range("L48").Value = 15
$ws.Range("D51").Value = 42
$ws.Range("E51").Value = 47
$ws.Range("H51").Value = 36
$ws.Range("G89").Value = 5
$ws.Range("C96").Value = 302
$ws.Range("D96").Value = 313
$ws.Range("E96").Value = 312
$ws.Range("F96").Value = 342
$ws.Range("G96").Value = 329
$ws.Range("H96").Value = 304
$ws.Range("I96").Value = 407
$ws.Range("J96").Value = 349
$ws.Range("K96").Value = 414
$ws.Range("L96").Value = 424

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range("D4").Value = 4
$ws.Range("J4").Value = 3
$ws.Range("D5").Value = 4
$ws.Range("J5").Value = 8

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("C2").Value = 3
$ws.Range("J3").Value = 5
$ws.Range("J5").Value = 12
$ws.Range("L5").Value = 10
$ws.Range("C6").Value = 24
$ws.Range("J6").Value = 25
$ws.Range("L6").Value = 14

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("F4").Value = 1
$ws.Range("J5").Value = 10
$ws.Range("F6").Value = 29
$ws.Range("J6").Value = 15

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("K2").Value = 6
$ws.Range("C6").Value = 13
$ws.Range("E6").Value = 10
$ws.Range("I6").Value = 12
$ws.Range("C7").Value = 15
$ws.Range("E7").Value = 13
$ws.Range("I7").Value = 22
$ws.Range("K7").Value = 29

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("L2").Value = 1
$ws.Range("L6").Value = 15

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("D6").Value = 23
$ws.Range("E6").Value = 36
$ws.Range("H6").Value = 26
$ws.Range("D7").Value = 42
$ws.Range("E7").Value = 47
$ws.Range("H7").Value = 36

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range("F6").Value = 4
$ws.Range("F7").Value = 5

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range("D4").Value = 3
$ws.Range("D5").Value = 4
